# "Switch pandoc to nix."
#
# The upstream diff for this commit is two kinds of change:
#
#  1. word/numbering.xml: the abstractNum(990) <w:nsid> value gets
#     zero-padded from "A990" to "0000A990" (a cosmetic artifact of the
#     docx writer used to regenerate the reference template under nix
#     vs. the previous toolchain). `nsid` is an internal numbering GUID
#     fragment that is not surfaced anywhere on the Word object model
#     (no Style/List/ListTemplate/ListLevel property maps to it), so it
#     is not reachable from COM automation.
#
#  2. word/styles.xml: a batch of "*Tok" character styles (pandoc's
#     syntax-highlighting styles) have their <w:rPr> child element order
#     normalized so boolean toggles (<w:b/>, <w:i/>) are serialized
#     before <w:color/>. Re-touching a style's Font.Bold / Font.Italic
#     through COM reproduces exactly this canonical re-serialization
#     order without altering the effective formatting (the values were
#     already bold/italic - only the on-disk element order changes).

$d = $word.ActiveDocument

# Styles whose <w:rPr> was just <w:color/><w:b/> -> now <w:b/><w:color/>
$boldOnly = @("KeywordTok", "ImportTok", "ControlFlowTok", "AlertTok", "ErrorTok")
foreach ($name in $boldOnly) {
    $d.Styles($name).Font.Bold = $true
}

# Styles whose <w:rPr> was just <w:color/><w:i/> -> now <w:i/><w:color/>
$italicOnly = @("CommentTok", "DocumentationTok")
foreach ($name in $italicOnly) {
    $d.Styles($name).Font.Italic = $true
}

# Styles whose <w:rPr> was <w:color/><w:b/><w:i/> -> now <w:b/><w:i/><w:color/>
$boldItalic = @("AnnotationTok", "CommentVarTok", "InformationTok", "WarningTok")
foreach ($name in $boldItalic) {
    $s = $d.Styles($name)
    $s.Font.Bold = $true
    $s.Font.Italic = $true
}
